# Applies the "error solve ifrs list" correction to the 두올 IFRS sheet:
# rows 2-6 get corrected (much smaller-scale) financial figures, with a
# handful of per-row columns now blank; rows 7-9 lose all of their data
# columns (D:AJ) entirely, keeping only the period label in column C.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("J2").ClearContents()
$ws.Range("Y2:Z2").ClearContents()
$ws.Range("AD2").ClearContents()
$ws.Range("AH2").ClearContents()
$ws.Range("D2").Value = 3168
$ws.Range("E2").Value = 184
$ws.Range("F2").Value = 184
$ws.Range("G2").Value = 123
$ws.Range("H2").Value = 124
$ws.Range("I2").Value = 124
$ws.Range("K2").Value = 2090
$ws.Range("L2").Value = 1648
$ws.Range("M2").Value = 442
$ws.Range("N2").Value = 443
$ws.Range("O2").Value = -1
$ws.Range("P2").Value = 279
$ws.Range("Q2").Value = 140
$ws.Range("R2").Value = -92
$ws.Range("S2").Value = -4
$ws.Range("T2").Value = 49
$ws.Range("U2").Value = 92
$ws.Range("V2").Value = 804
$ws.Range("W2").Value = 5.81
$ws.Range("X2").Value = 3.91
$ws.Range("AA2").Value = 372.74
$ws.Range("AB2").Value = 52.3
$ws.Range("AC2").Value = 661
$ws.Range("AE2").Value = 1944
$ws.Range("AF2").Value = 0
$ws.Range("AG2").Value = 0
$ws.Range("AI2").Value = 0
$ws.Range("AJ2").Value = 16587293

# Row 3
$ws.Range("AD3").ClearContents()
$ws.Range("AH3").ClearContents()
$ws.Range("D3").Value = 3663
$ws.Range("E3").Value = 261
$ws.Range("F3").Value = 261
$ws.Range("G3").Value = 234
$ws.Range("H3").Value = 209
$ws.Range("I3").Value = 206
$ws.Range("J3").Value = 3
$ws.Range("K3").Value = 2303
$ws.Range("L3").Value = 1416
$ws.Range("M3").Value = 887
$ws.Range("N3").Value = 885
$ws.Range("O3").Value = 1
$ws.Range("P3").Value = 383
$ws.Range("Q3").Value = 91
$ws.Range("R3").Value = -132
$ws.Range("S3").Value = 44
$ws.Range("T3").Value = 130
$ws.Range("U3").Value = -38
$ws.Range("V3").Value = 600
$ws.Range("W3").Value = 7.11
$ws.Range("X3").Value = 5.7
$ws.Range("Y3").Value = 31.03
$ws.Range("Z3").Value = 9.5
$ws.Range("AA3").Value = 159.74
$ws.Range("AB3").Value = 123.75
$ws.Range("AC3").Value = 904
$ws.Range("AE3").Value = 3883
$ws.Range("AF3").Value = 0
$ws.Range("AG3").Value = 0
$ws.Range("AI3").Value = 0
$ws.Range("AJ3").Value = 16587293

# Row 4
$ws.Range("D4").Value = 2959
$ws.Range("E4").Value = 184
$ws.Range("F4").Value = 184
$ws.Range("G4").Value = 179
$ws.Range("H4").Value = 127
$ws.Range("I4").Value = 123
$ws.Range("J4").Value = 4
$ws.Range("K4").Value = 2206
$ws.Range("L4").Value = 989
$ws.Range("M4").Value = 1217
$ws.Range("N4").Value = 1212
$ws.Range("O4").Value = 5
$ws.Range("P4").Value = 456
$ws.Range("Q4").Value = 115
$ws.Range("R4").Value = -75
$ws.Range("S4").Value = -19
$ws.Range("T4").Value = 79
$ws.Range("U4").Value = 36
$ws.Range("V4").Value = 362
$ws.Range("W4").Value = 6.21
$ws.Range("X4").Value = 4.3
$ws.Range("Y4").Value = 11.77
$ws.Range("Z4").Value = 5.65
$ws.Range("AA4").Value = 81.27
$ws.Range("AB4").Value = 166.65
$ws.Range("AC4").Value = 500
$ws.Range("AD4").Value = 6.24
$ws.Range("AE4").Value = 4575
$ws.Range("AF4").Value = 0.68
$ws.Range("AG4").Value = 77
$ws.Range("AH4").Value = 2.48
$ws.Range("AI4").Value = 16.46
$ws.Range("AJ4").Value = 26353521

# Row 5
$ws.Range("D5").Value = 3450
$ws.Range("E5").Value = 220
$ws.Range("F5").Value = 220
$ws.Range("G5").Value = 258
$ws.Range("H5").Value = 142
$ws.Range("I5").Value = 134
$ws.Range("J5").Value = 8
$ws.Range("K5").Value = 3320
$ws.Range("L5").Value = 1777
$ws.Range("M5").Value = 1542
$ws.Range("N5").Value = 1502
$ws.Range("O5").Value = 40
$ws.Range("P5").Value = 556
$ws.Range("Q5").Value = 613
$ws.Range("R5").Value = -476
$ws.Range("S5").Value = 54
$ws.Range("T5").Value = 159
$ws.Range("U5").Value = 454
$ws.Range("V5").Value = 575
$ws.Range("W5").Value = 6.39
$ws.Range("X5").Value = 4.12
$ws.Range("Y5").Value = 9.890000000000001
$ws.Range("Z5").Value = 5.14
$ws.Range("AA5").Value = 115.22
$ws.Range("AB5").Value = 186.73
$ws.Range("AC5").Value = 457
$ws.Range("AD5").Value = 8.91
$ws.Range("AE5").Value = 4624
$ws.Range("AF5").Value = 0.88
$ws.Range("AG5").Value = 67
$ws.Range("AH5").Value = 1.65
$ws.Range("AI5").Value = 16.13
$ws.Range("AJ5").Value = 27066680

# Row 6
$ws.Range("J6").ClearContents()
$ws.Range("O6").ClearContents()
$ws.Range("D6").Value = 4831
$ws.Range("E6").Value = 307
$ws.Range("F6").Value = 307
$ws.Range("G6").Value = 222
$ws.Range("H6").Value = 149
$ws.Range("I6").Value = 141
$ws.Range("K6").Value = 3466
$ws.Range("L6").Value = 1827
$ws.Range("M6").Value = 1640
$ws.Range("N6").Value = 1596
$ws.Range("P6").Value = 827
$ws.Range("Q6").Value = 390
$ws.Range("R6").Value = -96
$ws.Range("S6").Value = -239
$ws.Range("T6").Value = 115
$ws.Range("U6").Value = 274
$ws.Range("V6").Value = 494
$ws.Range("W6").Value = 6.35
$ws.Range("X6").Value = 3.08
$ws.Range("Y6").Value = 9.1
$ws.Range("Z6").Value = 4.38
$ws.Range("AA6").Value = 111.41
$ws.Range("AB6").Value = 101.86
$ws.Range("AC6").Value = 426
$ws.Range("AD6").Value = 6.71
$ws.Range("AE6").Value = 4912
$ws.Range("AF6").Value = 0.58
$ws.Range("AG6").Value = 80
$ws.Range("AH6").Value = 2.8
$ws.Range("AI6").Value = 18.44
$ws.Range("AJ6").Value = 27066680

# Row 7
$ws.Range("D7:AJ7").ClearContents()

# Row 8
$ws.Range("D8:AJ8").ClearContents()

# Row 9
$ws.Range("D9:AJ9").ClearContents()
